# Apply the "北京-漫展信息.xlsx" update:
#  1) Sheet "展览": bump several "想去人数" (F column) counters.
#  2) Sheet "演出": the first listed show (row 2) sold out / was removed from
#     the feed, so every later row's data (B:I) shifts up one row and the
#     now-duplicate last row (15) is removed. Column A (the static 0-based
#     index) is left untouched, which is why it is NOT part of the shift.
#  3) Sheet "全部类型": bump the same "想去人数" (F column) counters as sheet 1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) 展览 (Exhibitions) - F column counter bumps
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    3  = 32
    4  = 65
    5  = 85
    6  = 886
    7  = 476
    8  = 4792
    9  = 4792
    13 = 31
    15 = 136
    16 = 7729
    17 = 255
    20 = 542
    21 = 1423
    23 = 6293
    24 = 2263
    26 = 2095
    29 = 6218
    30 = 150
    31 = 37
    35 = 6553
    40 = 23
    41 = 37
    42 = 2476
    47 = 462
    48 = 2164
    49 = 52
    50 = 1096
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# ---------------------------------------------------------------------------
# 2) 演出 (Performances) - row 2 drops out of the feed; rows 3..15 (columns
#    B..I) shift up into rows 2..14, then the now-empty row 15 is deleted.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("B3:I15").Copy()
$ws2.Range("B2:I2").PasteSpecial()
$ws2.Application.CutCopyMode = $false
$ws2.Rows.Item(15).Delete()

# ---------------------------------------------------------------------------
# 3) 全部类型 (All types) - same F column counter bumps as 展览
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    4  = 32
    5  = 65
    7  = 85
    9  = 476
    10 = 4792
    11 = 4792
    15 = 31
    17 = 136
    18 = 7729
    19 = 7729
    20 = 255
    22 = 542
    23 = 1423
    25 = 6293
    26 = 2263
    27 = 2095
    29 = 6218
    30 = 150
    32 = 37
    36 = 6553
    41 = 23
    43 = 2476
    47 = 462
    49 = 2164
    50 = 52
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
